# Add four new English/Chinese word-pair rows to the end of the word list
# (Sheet1), then move/select the cell right after the newly added data,
# matching the upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New vocabulary pairs appended after the existing last data row (167).
$newPairs = @(
    @("microelectronics", "微电子"),
    @("parity", "平价"),
    @("frame", "框架"),
    @("sample", "样品")
)

$startRow = 168
for ($i = 0; $i -lt $newPairs.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newPairs[$i][0]
    $ws.Cells.Item($row, 2).Value = $newPairs[$i][1]
}

# Update the active selection to reflect where editing left off.
$ws.Range("C167").Select()
